$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" mc:Ignorable="w14" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006"><w:body><w:p w14:paraId="2C361703" w14:textId="321DA396" w:rsidR="00B80F4D" w:rsidRPr="00B80F4D" w:rsidRDefault="00B80F4D" w:rsidP="00B80F4D"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="40"/></w:rPr></w:pPr><w:r w:rsidRPr="00B80F4D"><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Rails</w:t></w:r></w:p><w:p w14:paraId="01866EC1" w14:textId="7CDD2080" w:rsidR="00456B02" w:rsidRDefault="00456B02"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Rails new &lt;name of app folder&gt;</w:t></w:r></w:p><w:p w14:paraId="1D539F7B" w14:textId="107D40F8" w:rsidR="00456B02" w:rsidRDefault="00456B02"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>-creates a new app folder in current pwd</w:t></w:r></w:p><w:p w14:paraId="78385C81" w14:textId="77777777" w:rsidR="00456B02" w:rsidRDefault="00456B02"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p w14:paraId="145D8541" w14:textId="6D444B2E" w:rsidR="005F6747" w:rsidRDefault="00C2329D"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00C2329D"><w:rPr><w:b/></w:rPr><w:t>Rails server</w:t></w:r></w:p><w:p w14:paraId="70051AEB" w14:textId="30B43BF1" w:rsidR="00C2329D" w:rsidRDefault="00C2329D"><w:pPr><w:rPr><w:rFonts w:ascii="Andale Mono" w:hAnsi="Andale Mono" w:cs="Andale Mono"/><w:color w:val="2FFF12"/></w:rPr></w:pPr><w:r><w:t>- starts running rails server</w:t></w:r><w:r w:rsidR="00664587"><w:t xml:space="preserve"> at </w:t></w:r><w:hyperlink r:id="rId4" w:history="1"><w:r w:rsidR="000520EC" w:rsidRPr="003C610B"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Andale Mono" w:hAnsi="Andale Mono" w:cs="Andale Mono"/></w:rPr><w:t>http://localhost:3000</w:t></w:r></w:hyperlink></w:p><w:p w14:paraId="18CC0E2F" w14:textId="77777777" w:rsidR="000520EC" w:rsidRDefault="000520EC"><w:pPr><w:rPr><w:rFonts w:ascii="Andale Mono" w:hAnsi="Andale Mono" w:cs="Andale Mono"/><w:color w:val="2FFF12"/></w:rPr></w:pPr></w:p><w:p w14:paraId="0D6C65B8" w14:textId="40CAB799" w:rsidR="000520EC" w:rsidRDefault="000520EC" w:rsidP="000520EC"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="006829E1"><w:rPr><w:b/></w:rPr><w:t>rail</w:t></w:r><w:r w:rsidR="006F6375"><w:rPr><w:b/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidRPr="006829E1"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> generate controller pages home</w:t></w:r></w:p><w:p w14:paraId="4D19C868" w14:textId="5A34B80F" w:rsidR="006829E1" w:rsidRDefault="006829E1" w:rsidP="000520EC"><w:r><w:t xml:space="preserve">- creates page at </w:t></w:r><w:hyperlink r:id="rId5" w:history="1"><w:r w:rsidR="00294D35" w:rsidRPr="003C610B"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>http://localhost:3000/pages/home</w:t></w:r></w:hyperlink></w:p><w:p w14:paraId="407160AA" w14:textId="2AE5F392" w:rsidR="00294D35" w:rsidRDefault="00294D35" w:rsidP="000520EC"><w:pPr><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">- is located on comp at </w:t></w:r><w:r w:rsidRPr="00294D35"><w:rPr><w:sz w:val="22"/></w:rPr><w:t>/Users/yehonatanmeschedekrasa/Desktop/pinterteresting</w:t></w:r><w:r w:rsidRPr="00294D35"><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>/app/views/pages</w:t></w:r></w:p><w:p w14:paraId="1A1A79BE" w14:textId="77777777" w:rsidR="00294D35" w:rsidRDefault="00294D35" w:rsidP="000520EC"><w:pPr><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="04BD1EBB" w14:textId="582FF7CD" w:rsidR="00294D35" w:rsidRDefault="006F6375" w:rsidP="000520EC"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>rails destroy controller pages home</w:t></w:r></w:p><w:p w14:paraId="3BF24438" w14:textId="43425D0C" w:rsidR="006F6375" w:rsidRDefault="006F6375" w:rsidP="000520EC"><w:r><w:t>-removes homepage created above</w:t></w:r></w:p><w:p w14:paraId="230EA3E8" w14:textId="2C794904" w:rsidR="0089532A" w:rsidRDefault="006F6375" w:rsidP="000520EC"><w:r><w:t>-don’t need to type “home”</w:t></w:r></w:p><w:p w14:paraId="2147D840" w14:textId="77777777" w:rsidR="00753456" w:rsidRDefault="00753456" w:rsidP="000520EC"/><w:p w14:paraId="19D2950F" w14:textId="77777777" w:rsidR="00753456" w:rsidRPr="00753456" w:rsidRDefault="00753456" w:rsidP="00753456"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00753456"><w:rPr><w:b/></w:rPr><w:t>rake routes</w:t></w:r></w:p><w:p w14:paraId="5C3FE21F" w14:textId="3EAD131A" w:rsidR="00753456" w:rsidRDefault="00753456" w:rsidP="00753456"><w:r><w:t>-shows available routes in current directory/app</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>rails generate scaffold pins</w:t></w:r></w:p><w:p><w:r><w:t>-adds pins to webapp</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Rails</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Http</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">create </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">    -</w:t></w:r><w:r><w:tab/><w:t>post</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= new + create</w:t></w:r></w:p><w:p><w:r><w:t>read</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">    -</w:t></w:r><w:r><w:tab/><w:t>get</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= show</w:t></w:r></w:p><w:p><w:r><w:t>update</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">    -</w:t></w:r><w:r><w:tab/><w:t>put</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= edit + update</w:t></w:r></w:p><w:p><w:r><w:t>destroy    -</w:t></w:r><w:r><w:tab/><w:t>delete</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>= destroy</w:t></w:r></w:p><w:p><w:r><w:t>(crud)</w:t></w:r></w:p><w:p><w:r><w:t>read all = index</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Heroku</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>heroku open</w:t></w:r></w:p><w:p><w:r><w:t>-</w:t></w:r><w:r><w:t>opens webapp at heroku url</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>heroku rename &lt; new name &gt;</w:t></w:r></w:p><w:p><w:r><w:t>-rename initial part of heroku url</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>git push</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> heroku master</w:t></w:r></w:p><w:p><w:r><w:t>-pushes to heroku app online</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">heroku logs </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>--</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>tail</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>-shows log files of your webapp</w:t></w:r></w:p><w:p><w:r><w:t>-will show you errors if you running into the error “We’re sorry but something went wrong”</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>heroku</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> run</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> rake db:migrate</w:t></w:r></w:p><w:p><w:r><w:t>-migrates database, aka updates database</w:t></w:r><w:r><w:t xml:space="preserve"> for heroku website</w:t></w:r></w:p><w:p w14:paraId="016DAB5F" w14:textId="2BAB5214" w:rsidR="00387ECE" w:rsidRDefault="00387ECE" w:rsidP="00753456"><w:r><w:t>-do if you’re getting a “Something went wrong message”</w:t></w:r></w:p><w:p w14:paraId="66789F95" w14:textId="3C8EFA8F" w:rsidR="00387ECE" w:rsidRPr="008C7886" w:rsidRDefault="00387ECE" w:rsidP="00753456"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>-need to reset server</w:t></w:r></w:p><w:p w14:paraId="5A0B4D74" w14:textId="77777777" w:rsidR="006A5211" w:rsidRDefault="006A5211" w:rsidP="00B80F4D"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="40"/></w:rPr></w:pPr></w:p><w:p w14:paraId="5E68E64D" w14:textId="1070FFA0" w:rsidR="00066CAE" w:rsidRDefault="006A5211" w:rsidP="00B80F4D"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Devise</w:t></w:r></w:p><w:p w14:paraId="103D5C7F" w14:textId="1D12291D" w:rsidR="000C67E7" w:rsidRPr="000C67E7" w:rsidRDefault="000C67E7" w:rsidP="00B80F4D"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>12, 13, 14</w:t></w:r></w:p><w:p w14:paraId="410AEC74" w14:textId="779B5CC1" w:rsidR="000C67E7" w:rsidRDefault="000C67E7" w:rsidP="000C67E7"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>rails generate devise user</w:t></w:r></w:p><w:p w14:paraId="5B1DF4AE" w14:textId="73CCBB61" w:rsidR="000C67E7" w:rsidRDefault="000C67E7" w:rsidP="000C67E7"><w:r><w:t>-</w:t></w:r><w:r><w:t>creates a user</w:t></w:r></w:p><w:p w14:paraId="7934FB43" w14:textId="77777777" w:rsidR="00C66ADD" w:rsidRDefault="00C66ADD" w:rsidP="000C67E7"/><w:p w14:paraId="28C7A441" w14:textId="532BDDD8" w:rsidR="000C67E7" w:rsidRDefault="000C67E7" w:rsidP="000C67E7"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>rails destroy</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> devise user</w:t></w:r></w:p><w:p w14:paraId="0DF5277C" w14:textId="276062C0" w:rsidR="000C67E7" w:rsidRDefault="000C67E7" w:rsidP="000C67E7"><w:r><w:t>-</w:t></w:r><w:r><w:t>deletes</w:t></w:r><w:r><w:t xml:space="preserve"> a user</w:t></w:r></w:p><w:p w14:paraId="5AD7EDB8" w14:textId="77777777" w:rsidR="000C67E7" w:rsidRDefault="000C67E7" w:rsidP="000C67E7"/><w:p w14:paraId="598103E1" w14:textId="2E451A66" w:rsidR="00352DFD" w:rsidRDefault="00352DFD" w:rsidP="00352DFD"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>rake db:migrate</w:t></w:r></w:p><w:p w14:paraId="457EBADC" w14:textId="617A7B31" w:rsidR="00352DFD" w:rsidRDefault="00352DFD" w:rsidP="00352DFD"><w:r><w:t>-</w:t></w:r><w:r><w:t>migrates database, aka updates database</w:t></w:r></w:p><w:p w14:paraId="586BC109" w14:textId="77777777" w:rsidR="00EE1FE9" w:rsidRDefault="00EE1FE9" w:rsidP="00352DFD"/><w:p w14:paraId="0863BFFD" w14:textId="3BFE2E6D" w:rsidR="00EE1FE9" w:rsidRDefault="00EE1FE9" w:rsidP="00352DFD"/><w:p w14:paraId="74C2353B" w14:textId="77777777" w:rsidR="00352DFD" w:rsidRPr="00553832" w:rsidRDefault="00352DFD" w:rsidP="000C67E7"/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
Write-Output $d.Paragraphs.Count
